$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments ---
$ws.Columns.Item(3).ColumnWidth = 65.42578125
$ws.Columns.Item(6).ColumnWidth = 53.7109375

# --- Bug report rows (text / numeric values) ---
# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "high"
$ws.Range("C2").Value = "error page on po creation"
$ws.Range("D2").Value = "purchaseordercontroler.cs"
$ws.Range("F2").Value = "create po and submit it then try to create a second one"
$ws.Range("H2").Value = "unfixed"
$ws.Range("I2").Value = "Mathew"

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "medium"
$ws.Range("C3").Value = "no validation on po creation"
$ws.Range("D3").Value = "purchaseordercontroler.cs"
$ws.Range("F3").Value = "create po with none valid input"
$ws.Range("H3").Value = "unfixed"
$ws.Range("I3").Value = "Mathew"

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "high"
$ws.Range("C4").Value = "page for /RECEIVING_LOG/findPO/orig not found when looking up a rl"
$ws.Range("D4").Value = "receivinglogcontroler.cs"
$ws.Range("F4").Value = "look up any record in receiving log"
$ws.Range("H4").Value = "unfixed"
$ws.Range("I4").Value = "Mathew"

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "medium"
$ws.Range("C5").Value = "no validation on begin receiving log"
$ws.Range("D5").Value = "receivinglogcontroler.cs"
$ws.Range("F5").Value = "look up any record with incorrect information"
$ws.Range("H5").Value = "unfixed"
$ws.Range("I5").Value = "Mathew"

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "low"
$ws.Range("C6").Value = "deletes invalid input in the MSRP field and states it is required"
$ws.Range("D6").Value = "itemmanagementcontroler.cs"
$ws.Range("F6").Value = "input 'something' into MSRP on item creation"
$ws.Range("H6").Value = "unfixed"
$ws.Range("I6").Value = "Nate"

# --- Reported On dates (column J) ---
# Set the format + value on J2 once, then propagate the same style to J3:J6
# via copy/paste-format so that only a single new style entry is created
# (matching the single extra cellXfs entry produced by the real edit).
$dateMid = Get-Date -Year 2014 -Month 2 -Day 21 -Hour 0 -Minute 0 -Second 0
$dateLow = Get-Date -Year 2014 -Month 2 -Day 17 -Hour 0 -Minute 0 -Second 0

$j2 = $ws.Range("J2")
$j2.NumberFormat = "mm-dd-yy"
$j2.Value = $dateMid

$j2.Copy()
$ws.Range("J3:J6").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("J3").Value = $dateMid
$ws.Range("J4").Value = $dateMid
$ws.Range("J5").Value = $dateMid
$ws.Range("J6").Value = $dateLow
